$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update label text: "Hours at the beginning of sprint of 1" -> "... of 2" (cell L24) ---
$ws.Range("L24").Value = "Hours at the beginning of sprint of 2"

# --- Row 4: "Effort remaining" data ---
$ws.Range("F4").Value = 55
$ws.Range("G4").Value = 53
$ws.Range("H4").Value = 53
$ws.Range("I4").Value = 49
$ws.Range("J4").Value = 45
$ws.Range("K4").Value = 35
$ws.Range("L4").Value = 31
$ws.Range("M4").Value = 26
$ws.Range("N4").Value = 16
$ws.Range("O4").Value = 11
$ws.Range("P4").Value = 5
$ws.Range("Q4").Value = 5
$ws.Range("R4").Value = 5

# --- Row 18: "Estimate Effort" data ---
$ws.Range("C18").Value = 59
$ws.Range("D18").Value = 59
$ws.Range("G18").Value = 53
$ws.Range("H18").Value = 53
$ws.Range("I18").Value = 48
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = 37
$ws.Range("L18").Value = 31
$ws.Range("M18").Value = 28
$ws.Range("N18").Value = 16
$ws.Range("O18").Value = 13
$ws.Range("P18").Value = 5
$ws.Range("Q18").Value = 5
$ws.Range("R18").Value = 5

# --- Row 19: "Actual Effort" data ---
$ws.Range("E19").Value = 57
$ws.Range("F19").Value = 55
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = 53
$ws.Range("I19").Value = 49
$ws.Range("J19").Value = 45
$ws.Range("K19").Value = 35
$ws.Range("L19").Value = 31
$ws.Range("M19").Value = 26
$ws.Range("N19").Value = 16
$ws.Range("O19").Value = 11
$ws.Range("P19").Value = 5
$ws.Range("Q19").Value = 0
$ws.Range("R19").Value = 0

# --- Row 25: single value update ---
$ws.Range("K25").Value = 5

# --- Sheet view: scroll position and selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("Q43").Select()
